$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Row 2 is a new data row appended below the header row: it repeats
# the same five header values (name / bank+branch / account / date /
# type) plus a leading sequence number in column A.
# ------------------------------------------------------------------

# --- Scratch step: manufacture a "plain" cell style (a style record
# that is visually identical to the default style but lives at its
# own distinct index) by touching an alignment sub-property and then
# reverting it. The scratch cell/row is removed afterwards so it
# leaves no trace in the sheet, but the style definition itself
# survives in the workbook's style table for reuse below.
$scratch = $ws.Range("Z100")
$scratch.Value = 0
$scratch.WrapText = $false
$scratch.Copy()
$ws.Range("B2:G2").PasteSpecial(-4122)
$scratch.EntireRow.Delete()

# --- A2: sequence number, reusing the header's bordered/centered style.
$ws.Range("A2").Value = 13
$ws.Range("C1").Copy()
$ws.Range("A2").PasteSpecial(-4122)

# --- B2:G2: same values as B1:G1 (plain style already applied above).
$ws.Range("B2").Value = $ws.Range("B1").Value2
$ws.Range("C2").Value = "'" + $ws.Range("C1").Value2
$ws.Range("D2").Value = "'" + $ws.Range("D1").Value2
$ws.Range("E2").Value = "'" + $ws.Range("E1").Value2
$ws.Range("F2").Value = "'" + $ws.Range("F1").Value2
$ws.Range("G2").Value = "'" + $ws.Range("G1").Value2
